$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data
# Each target cell is forced to Text format ("@") before assignment so that
# numeric-looking strings (e.g. "1.000", "0.9999") are preserved as text,
# matching the original inline-string cell contents.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.557.72'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.68%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.877.49'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.01%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.19%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.51'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.79%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.20%  '
# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.06%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2917'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.50%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06505'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.10%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.98'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.71%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07745'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.39%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7391'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.06%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '96.54'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.37%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.873.36'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.21%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.185'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.21%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '274.05'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.89%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.649.47'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.07%  '
# Row 18
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.01%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9999'
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007518'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.08%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.119.49'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.56%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.28%  '
# Row 23
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.45%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.199'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.71%  '
# Row 25
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '165.37'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.59%  '
# Row 26
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.194'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.20%  '
# Row 27
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.58%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.909'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.63%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.09843'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.93%  '
# Row 30
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.36%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.502'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.78%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.270'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.55%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.099'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.47%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04815'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.56%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.124'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.29%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6957'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.04%  '
# Row 37
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.23%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01868'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.95%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.763'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.50%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.274'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.47%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.56'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.27%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.986'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.35%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4217'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.45%  '
# Row 44
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.18%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8350'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.02%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.79'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.19%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.418'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.72%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.988'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.17%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.32'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.65%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '917.08'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.95%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05676'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.71%  '
